## Update gh-pages output data (generated at 456a3b4)
## Applies the "想去人数" (F) / "最低票价" (G) refresh, plus one swapped-in
## activity (row 45) on sheet "展览", across all four worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 1886
$ws1.Range("F3").Value  = 26
$ws1.Range("F5").Value  = 762
$ws1.Range("F7").Value  = 507
$ws1.Range("F8").Value  = 920
$ws1.Range("F9").Value  = 1617
$ws1.Range("F10").Value = 1285
$ws1.Range("F11").Value = 1550
$ws1.Range("F12").Value = 68
$ws1.Range("F13").Value = 1547
$ws1.Range("F14").Value = 346
$ws1.Range("F15").Value = 1694
$ws1.Range("F17").Value = 1125
$ws1.Range("F18").Value = 377
$ws1.Range("F20").Value = 112
$ws1.Range("F21").Value = 1796
$ws1.Range("F22").Value = 248
$ws1.Range("F26").Value = 1251
$ws1.Range("F27").Value = 1067
$ws1.Range("F29").Value = 577
$ws1.Range("F30").Value = 1171
$ws1.Range("F33").Value = 1173
$ws1.Range("F34").Value = 1119
$ws1.Range("F35").Value = 284
$ws1.Range("F36").Value = 84
$ws1.Range("F37").Value = 886
$ws1.Range("F38").Value = 1694
$ws1.Range("F42").Value = 2060
$ws1.Range("F44").Value = 837

# Row 45 is a brand-new entry replacing the old one in the source feed
$ws1.Range("C45").Value = "上海·创造力动漫游戏嘉年华1.0"
$ws1.Range("D45").Value = "莘福路288号 美莘商业广场"
$ws1.Range("F45").Value = 2
$ws1.Range("G45").Value = 65
$ws1.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=87667"
$ws1.Range("I45").Value = "//i2.hdslb.com/bfs/openplatform/202406/cjmOiK0E1718378936182.png"

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value  = 51
$ws2.Range("G5").Value  = 380
$ws2.Range("F6").Value  = 1502
$ws2.Range("F9").Value  = 2601
$ws2.Range("F10").Value = 1223
$ws2.Range("F20").Value = 24
$ws2.Range("F22").Value = 322
$ws2.Range("F23").Value = 88480
$ws2.Range("F24").Value = 31
$ws2.Range("F29").Value = 253
$ws2.Range("F31").Value = 226
$ws2.Range("F35").Value = 19
$ws2.Range("F44").Value = 142

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F5").Value  = 2917
$ws3.Range("F6").Value  = 4685
$ws3.Range("F7").Value  = 145
$ws3.Range("F9").Value  = 590
$ws3.Range("F10").Value = 768
$ws3.Range("F11").Value = 480
$ws3.Range("F12").Value = 407
$ws3.Range("F13").Value = 1140
$ws3.Range("F14").Value = 312
$ws3.Range("F15").Value = 742

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 1886
$ws4.Range("F4").Value  = 26
$ws4.Range("F5").Value  = 4685
$ws4.Range("F6").Value  = 768
$ws4.Range("F8").Value  = 407
$ws4.Range("F9").Value  = 407
$ws4.Range("F10").Value = 507
$ws4.Range("F11").Value = 920
$ws4.Range("F12").Value = 1223
$ws4.Range("F13").Value = 1617
$ws4.Range("F14").Value = 1285
$ws4.Range("F15").Value = 1550
$ws4.Range("F16").Value = 68
$ws4.Range("F17").Value = 1547
$ws4.Range("F20").Value = 1695
$ws4.Range("F21").Value = 1125
$ws4.Range("F22").Value = 377
$ws4.Range("F24").Value = 742
$ws4.Range("F25").Value = 742
$ws4.Range("F26").Value = 1796
$ws4.Range("F27").Value = 248
$ws4.Range("F31").Value = 1251
$ws4.Range("F32").Value = 322
$ws4.Range("F33").Value = 1067
$ws4.Range("F35").Value = 1171
$ws4.Range("F37").Value = 1173
$ws4.Range("F39").Value = 1119
$ws4.Range("F40").Value = 284
$ws4.Range("F41").Value = 886
$ws4.Range("F43").Value = 1694
$ws4.Range("F46").Value = 2060
$ws4.Range("F48").Value = 837
